# Adding the Bernoulli distribution row to the dist_table sheet.
#
# The new row is inserted as row 2 (pushing every existing row down by
# one), and an autofilter (with its matching _xlnm._FilterDatabase
# defined name) is turned on over the whole table, matching the
# upstream commit "Adding the Bernoulli distribution".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new row -----------------------------------
$ws.Rows.Item(2).Insert()

# Row 2 inherited row 1's (header) bold font + formatting from the
# insert; reset it back to the plain body-row look used by every other
# data row (General for the text columns, "@" text format for the two
# TRUE/FALSE flag columns).
$ws.Range("A2:H2").Font.Bold = $false
$ws.Range("A2:F2").NumberFormat = "General"
$ws.Range("G2:H2").NumberFormat = "@"

# --- 2. Fill in the Bernoulli row data --------------------------------
$ws.Range("A2").Value = "Bernoulli"
$ws.Range("B2").Value = "Standard"
$ws.Range("C2").Value = "bernoulli"
$ws.Range("D2").Value = "std"
$ws.Range("E2").Value = "binary"
$ws.Range("F2").Value = "uni"

# G2/H2 both need the literal text "TRUE" (not the boolean TRUE) to
# match the rest of the sheet, which stores these flags as text.
# Copying from an existing TRUE/TRUE row (row 9, "Geometric"/"Mean")
# carries over both the literal text value and its cell style exactly.
$ws.Range("G9:H9").Copy($ws.Range("G2:H2"))

# --- 3. Turn on the autofilter over the whole (now 21-row) table -----
$ws.Range("A1:H21").AutoFilter()
$ws.Names.Add("_xlnm._FilterDatabase", "=dist_table!`$A`$1:`$H`$21")

# --- 4. Misc view state that shifted along with the new row ----------
$ws.Range("B27").Select()
